# Doing Updates for Financials
# Update the OXFD yearly financials with the latest reported figures
# (column D = most recent period) and mark a handful of cells that no
# longer have a reported value as "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Income Statement -------------------------------------------------
$ws.Range("D8").Value  = 54700    # Total Revenue
$ws.Range("D9").Value  = 18500    # Cost of Revenue
$ws.Range("D10").Value = 36300    # Gross Profit
$ws.Range("D12").Value = 27500    # Research Development
$ws.Range("D17").Value = 81200    # Total Operating Expenses
$ws.Range("D18").Value = -26400   # Operating Income or Loss
$ws.Range("D21").Value = -27400   # Earnings Before Interest And Taxes
$ws.Range("J21").Value = "NA"     # Earnings Before Interest And Taxes (oldest period)
$ws.Range("D23").Value = -31600   # Income Before Tax
$ws.Range("D24").Value = 1600     # Income Tax Expense
$ws.Range("D26").Value = -33200   # Income After Tax
$ws.Range("D27").Value = -33200   # Net Income From Continuing Ops
$ws.Range("D29").Value = -21100   # Discontinued Operations
$ws.Range("D33").Value = -54300   # Net Income
$ws.Range("D35").Value = -54300   # Net Income Applicable To Common Shares

# --- Balance Sheet ------------------------------------------------------
$ws.Range("D43").Value = 23000    # Net Receivables
$ws.Range("D44").Value = 7100     # Inventory
$ws.Range("D45").Value = 17000    # Other Current Assets
$ws.Range("D48").Value = 11800    # Property Plant and Equipment
$ws.Range("D49").Value = 11400    # Goodwill
$ws.Range("D52").Value = 17700    # Other Assets
$ws.Range("D57").Value = 5600     # Accounts Payable
$ws.Range("D59").Value = 17500    # Other Current Liabilities

# --- Cash Flow Statement -------------------------------------------------
$ws.Range("D81").Value = -54300   # Net Income
$ws.Range("J83").Value  = "NA"    # Depreciation (oldest period)
$ws.Range("J94").Value  = "NA"    # Total Cash Flows From Investing Activities (oldest period)
$ws.Range("J100").Value = "NA"    # Total Cash Flows From Financing Activities (oldest period)
$ws.Range("J101").Value = "NA"    # Effect Of Exchange Rate Changes (oldest period)
